# Fruta / hortaliza, semanal
# Insert the latest week's two records (Primera / Segunda) at the top of the
# data block for this product/market, pushing the historical rows down by
# two rows (dimension grows from A1:R528 to A1:R530).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 494; everything from 494:528 shifts
# down to 496:530 and the used range / dimension is updated automatically.
$ws.Range("A494:A495").EntireRow.Insert()

# New row 494 - "Primera" grade record for 2022-06-02.
$ws.Range("A494").Value = 9
$ws.Range("B494").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C494").Value = "Metropolitana"
$ws.Range("D494").Value = 44714
$ws.Range("E494").Value = 13
$ws.Range("F494").Value = 100114014
$ws.Range("G494").Value = "Betarraga"
$ws.Range("H494").Value = "Sin especificar"
$ws.Range("I494").Value = "Primera"
$ws.Range("J494").Value = 10600
$ws.Range("K494").Value = 110
$ws.Range("L494").Value = 120
$ws.Range("M494").Value = 115
$ws.Range("N494").Value = "$/unidad"
$ws.Range("O494").Value = "Región Metropolitana"
$ws.Range("P494").Value = 115
$ws.Range("Q494").Value = 1
$ws.Range("R494").Value = "Hortaliza"

# New row 495 - "Segunda" grade record for 2022-06-02.
$ws.Range("A495").Value = 9
$ws.Range("B495").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C495").Value = "Metropolitana"
$ws.Range("D495").Value = 44714
$ws.Range("E495").Value = 13
$ws.Range("F495").Value = 100114014
$ws.Range("G495").Value = "Betarraga"
$ws.Range("H495").Value = "Sin especificar"
$ws.Range("I495").Value = "Segunda"
$ws.Range("J495").Value = 4300
$ws.Range("K495").Value = 90
$ws.Range("L495").Value = 90
$ws.Range("M495").Value = 90
$ws.Range("N495").Value = "$/unidad"
$ws.Range("O495").Value = "Región Metropolitana"
$ws.Range("P495").Value = 90
$ws.Range("Q495").Value = 1
$ws.Range("R495").Value = "Hortaliza"
